# Add 2022-Q3 data.
#
# Before:  总计 , 2022-Q1
# After:   总计 , 2022-Q3 , 2022-Q1
#
# The existing "2022-Q1" detail sheet becomes the new "2022-Q3" sheet (same
# slot, new content), and a brand-new sheet named "2022-Q1" is appended right
# after it, re-populated with the fund-holding detail data that used to live
# on the original "2022-Q1" sheet.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Capture the old "2022-Q1" detail-sheet data before we overwrite it.
# ---------------------------------------------------------------------------
$oldQ1 = $wb.Worksheets.Item("2022-Q1")

$oldRows = @(
    ,@("012348", "天弘恒生科技指数型发起式证券投资基金（QDII）A", "38.10", "92.34", "7.27", "2.7699", 4)
    ,@("968029", "恒生指数基金M类人民币（对冲）份额", "25.09", "97.94", "7.98", "2.0022", 2)
    ,@("012349", "天弘恒生科技指数型发起式证券投资基金（QDII）C", "14.77", "92.34", "7.27", "1.0738", 4)
    ,@("009562", "工银瑞信中国机会全球配置股票(QDII)美元", "6.65", "92.85", "7.01", "0.4662", 1)
    ,@("486001", "工银瑞信中国机会全球配置股票(QDII)", "6.65", "92.85", "7.01", "0.4662", 1)
    ,@("009563", "工银瑞信中国机会全球配置股票(QDII)港币", "6.65", "92.85", "7.01", "0.4662", 1)
    ,@("009225", "天弘中证中美互联网指数（QDII）A", "1.84", "94.90", "9.62", "0.1770", 3)
    ,@("009226", "天弘中证中美互联网指数（QDII）C", "0.59", "94.90", "9.62", "0.0568", 3)
)

# ---------------------------------------------------------------------------
# 2. Insert the brand new sheet right after "2022-Q1", then rename the two
#    sheets so the new slot carries the "2022-Q1" name and the existing slot
#    becomes "2022-Q3" - matching the final tab order/sheet identities.
# ---------------------------------------------------------------------------
$newQ1 = $wb.Worksheets.Add($null, $oldQ1)
$oldQ1.Name = "2022-Q3"
$newQ1.Name = "2022-Q1"

# ---------------------------------------------------------------------------
# 3. Re-write the new "2022-Q1" sheet with the preserved original data.
# ---------------------------------------------------------------------------
$newQ1.Range("B1").Value = "基金代码"
$newQ1.Range("C1").Value = "基金名称"
$newQ1.Range("D1").Value = "基金规模"
$newQ1.Range("E1").Value = "股票总仓位"
$newQ1.Range("F1").Value = "仓位占比"
$newQ1.Range("G1").Value = "持有市值(亿元)"
$newQ1.Range("H1").Value = "仓位排名"
$newQ1.Range("B1:H1").Style = "Bold_Border"

for ($i = 0; $i -lt $oldRows.Count; $i++) {
    $r = $i + 2
    $row = $oldRows[$i]
    $newQ1.Cells.Item($r, 1).Value = $i
    $newQ1.Cells.Item($r, 2).Value = $row[0]
    $newQ1.Cells.Item($r, 3).Value = $row[1]
    $newQ1.Cells.Item($r, 4).Value = $row[2]
    $newQ1.Cells.Item($r, 5).Value = $row[3]
    $newQ1.Cells.Item($r, 6).Value = $row[4]
    $newQ1.Cells.Item($r, 7).Value = $row[5]
    $newQ1.Cells.Item($r, 8).Value = $row[6]
}

# ---------------------------------------------------------------------------
# 4. Overwrite the (renamed) "2022-Q3" sheet with the new quarter's data.
# ---------------------------------------------------------------------------
$q3 = $wb.Worksheets.Item("2022-Q3")

$q3Rows = @(
    ,@("012349", "天弘恒生科技指数（QDII）C", "33.57", "92.84", "7.75", "2.6017", 4)
    ,@("012348", "天弘恒生科技指数（QDII）A", "30.64", "92.84", "7.75", "2.3746", 4)
    ,@("968029", "恒生指数基金M类人民币（对冲）份额", "27.03", "99.07", "7.61", "2.0570", 3)
    ,@("009562", "工银全球股票（QDII）美元", "5.89", "93.72", "6.25", "0.3681", 1)
    ,@("009563", "工银全球股票（QDII）港币", "5.89", "93.72", "6.25", "0.3681", 1)
    ,@("486001", "工银瑞信中国机会全球配置股票（QDII）人民币", "5.89", "93.72", "6.25", "0.3681", 1)
    ,@("009225", "天弘中证中美互联网指数（QDII）A", "1.20", "94.98", "8.25", "0.0990", 4)
    ,@("002379", "工银瑞信香港中小盘股票（QDII）人民币", "1.58", "78.58", "5.84", "0.0923", 2)
    ,@("002380", "工银瑞信香港中小盘股票（QDII）美元", "1.58", "78.58", "5.84", "0.0923", 2)
    ,@("009226", "天弘中证中美互联网指数（QDII）C", "0.60", "94.98", "8.25", "0.0495", 4)
)

for ($i = 0; $i -lt $q3Rows.Count; $i++) {
    $r = $i + 2
    $row = $q3Rows[$i]
    $q3.Cells.Item($r, 1).Value = $i
    $q3.Cells.Item($r, 2).Value = $row[0]
    $q3.Cells.Item($r, 3).Value = $row[1]
    $q3.Cells.Item($r, 4).Value = $row[2]
    $q3.Cells.Item($r, 5).Value = $row[3]
    $q3.Cells.Item($r, 6).Value = $row[4]
    $q3.Cells.Item($r, 7).Value = $row[5]
    $q3.Cells.Item($r, 8).Value = $row[6]
}

# ---------------------------------------------------------------------------
# 5. Update the "总计" (summary) sheet: insert a new row for 2022-Q3 right
#    above the existing 2022-Q1 row, and bump counts accordingly.
# ---------------------------------------------------------------------------
$summary = $wb.Worksheets.Item("总计")
$summary.Rows.Item(2).Insert()

$summary.Cells.Item(2, 1).Value = 0
$summary.Cells.Item(2, 2).Value = "2022-Q3"
$summary.Cells.Item(2, 3).Value = 10
$summary.Cells.Item(2, 4).Value = 8.47

$summary.Cells.Item(3, 1).Value = 1
$summary.Cells.Item(3, 2).Value = "2022-Q1"
